$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '29.549.25'
$ws.Range("E2").Value = '  +1.96%  '

$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '1.842.39'
$ws.Range("E3").Value = '  +0.81%  '

$ws.Range("D4").NumberFormat = '@'
$ws.Range("D4").Value = '0.9988'
$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("E5").Value = '  +1.06%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '0.6291'
$ws.Range("E6").Value = '  +2.14%  '

$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.9999'
$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.07443'
$ws.Range("E8").Value = '  +0.87%  '

$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.2954'
$ws.Range("E9").Value = '  +0.76%  '

$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '23.52'
$ws.Range("E10").Value = '  +2.59%  '

$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.07681'
$ws.Range("E11").Value = '  +0.51%  '

$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '1.848.37'
$ws.Range("E12").Value = '  +0.27%  '

$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '5.035'
$ws.Range("E13").Value = '  +1.14%  '

$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '0.6802'
$ws.Range("E14").Value = '  +1.59%  '

$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '83.55'
$ws.Range("E15").Value = '  +1.39%  '

$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '0.000009191'
$ws.Range("E16").Value = '  +1.96%  '

$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '5.939'
$ws.Range("E17").Value = '  +1.02%  '

$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '29.521.37'
$ws.Range("E18").Value = '  +1.83%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '2.100.96'
$ws.Range("E19").Value = '  +0.73%  '

$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '247.39'
$ws.Range("E20").Value = '  +5.01%  '

$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '12.60'
$ws.Range("E21").Value = '  -0.42%  '

$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '0.9999'
$ws.Range("E22").Value = '  +0.06%  '

$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '7.433'
$ws.Range("E23").Value = '  +3.68%  '

$ws.Range("E24").Value = '  +0.21%  '

$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '158.67'
$ws.Range("E25").Value = '  +0.23%  '

$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '0.1417'
$ws.Range("E26").Value = '  -0.33%  '

$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '8.595'
$ws.Range("E27").Value = '  +1.61%  '

$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '17.87'
$ws.Range("E28").Value = '  +0.59%  '

$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '1.499'
$ws.Range("E29").Value = '  +0.50%  '

$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '0.05954'
$ws.Range("E30").Value = '  +7.00%  '

$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '1.235'
$ws.Range("E31").Value = '  +2.04%  '

$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '4.120'
$ws.Range("E32").Value = '  +0.65%  '

$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '4.136'
$ws.Range("E33").Value = '  +0.62%  '

$ws.Range("E34").Value = '  +1.81%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '1.148'
$ws.Range("E35").Value = '  +1.35%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '0.7294'
$ws.Range("E36").Value = '  -1.47%  '

$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '2.616'
$ws.Range("E37").Value = '  -0.95%  '

$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '2.891'
$ws.Range("E38").Value = '  +4.19%  '

$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '1.233.84'
$ws.Range("E39").Value = '  +2.34%  '

$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '0.01768'
$ws.Range("E40").Value = '  -0.16%  '

$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '6.340'
$ws.Range("E41").Value = '  -0.08%  '

$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '0.9173'
$ws.Range("E42").Value = '  +2.32%  '

$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.38%  '

$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '2.014.24'
$ws.Range("E44").Value = '  +1.55%  '

$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '102.08'
$ws.Range("E45").Value = '  +1.02%  '

$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '65.96'
$ws.Range("E46").Value = '  +1.59%  '

$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '0.5063'
$ws.Range("E47").Value = '  -0.32%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '9.257'
$ws.Range("E48").Value = '  +1.64%  '

$ws.Range("B49").Value = 'TheSandbox'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '0.4074'
$ws.Range("E49").Value = '  +0.86%  '

$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '0.00000000118'
$ws.Range("E50").Value = '  -3.62%  '

$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '0.1166'
$ws.Range("E51").Value = '  +6.22%  '
